# Update cryptocurrency price/volume table (columns D and E) for rows 2-51
# Cells whose new value is a plain decimal number need NumberFormat "@"
# forced first so Excel stores them as text (matching source data), not a float.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.031.96"
$ws.Range("E2").Value = "  -1.54%  "
$ws.Range("D3").Value = "1.565.78"
$ws.Range("E3").Value = "  -0.38%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "206.65"
$ws.Range("E5").Value = "  -0.40%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.488"
$ws.Range("E6").Value = "  -1.53%  "
$ws.Range("E7").Value = "  +0.22%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.15"
$ws.Range("E8").Value = "  +0.31%  "
$ws.Range("E9").Value = "  -0.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0585"
$ws.Range("E10").Value = "  -0.59%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0867"
$ws.Range("E11").Value = "  +0.36%  "
$ws.Range("D12").Value = "1.792.29"
$ws.Range("E12").Value = "  -0.26%  "
$ws.Range("D13").Value = "1.571.56"
$ws.Range("E13").Value = "  -0.38%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "3.75"
$ws.Range("E14").Value = "  -1.79%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.514"
$ws.Range("E15").Value = "  -1.44%  "
$ws.Range("D16").Value = "27.097.87"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "61.78"
$ws.Range("E17").Value = "  -2.16%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "7.34"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "213.34"
$ws.Range("E19").Value = "  -0.24%  "
$ws.Range("D20").Value = "0.0₃0681"
$ws.Range("E20").Value = "  -1.17%  "
$ws.Range("E21").Value = "  +0.07%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.10"
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.38"
$ws.Range("E23").Value = "  -3.98%  "
$ws.Range("E24").Value = "  +0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "152.26"
$ws.Range("E25").Value = "  -0.47%  "
$ws.Range("E26").Value = "  -1.86%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "14.86"
$ws.Range("E27").Value = "  -1.24%  "
$ws.Range("E28").Value = "  +0.06%  "
$ws.Range("E29").Value = "  -1.64%  "
$ws.Range("E30").Value = "  -3.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0460"
$ws.Range("E31").Value = "  -2.23%  "
$ws.Range("D33").Value = "1.392.50"
$ws.Range("E33").Value = "  +2.50%  "
$ws.Range("E34").Value = "  -0.95%  "
$ws.Range("E35").Value = "  +0.18%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.30"
$ws.Range("E36").Value = "  -0.70%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.944"
$ws.Range("E37").Value = "  -2.71%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.0163"
$ws.Range("E38").Value = "  -2.22%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.808"
$ws.Range("E39").Value = "  -1.99%  "
$ws.Range("E40").Value = "  -3.66%  "
$ws.Range("E41").Value = "  +0.10%  "
$ws.Range("E42").Value = "  +3.25%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.42"
$ws.Range("E43").Value = "  +1.93%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.80"
$ws.Range("E44").Value = "  +0.70%  "
$ws.Range("E45").Value = "  +0.92%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "63.54"
$ws.Range("E46").Value = "  -1.10%  "
$ws.Range("D47").Value = "1.704.18"
$ws.Range("E47").Value = "  -0.41%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "85.30"
$ws.Range("E48").Value = "  -0.09%  "
$ws.Range("D49").Value = "0.0₇0974"
$ws.Range("E49").Value = "  -2.09%  "
$ws.Range("E50").Value = "  +0.22%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0948"
$ws.Range("E51").Value = "  -0.83%  "
